$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (diff-derived updates to the cryptos price table)
$updates = @{
    'D2' = '42.986.27'
    'E2' = '  -0.91%  '
    'D3' = '2.336.05'
    'E3' = '  +1.28%  '
    'E4' = '  +0.00%  '
    'D5' = '306.50'
    'E5' = '  -1.55%  '
    'D6' = '100.92'
    'E6' = '  -1.76%  '
    'D7' = '0.513'
    'E7' = '  -3.83%  '
    'E8' = '  +0.01%  '
    'D9' = '0.510'
    'E9' = '  -3.69%  '
    'D10' = '34.79'
    'E10' = '  -2.74%  '
    'D11' = '52.27'
    'E11' = '  +1.40%  '
    'D12' = '0.0800'
    'E12' = '  -1.70%  '
    'E13' = '  +0.63%  '
    'E14' = '  -2.30%  '
    'D15' = '15.83'
    'E15' = '  +5.41%  '
    'D16' = '2.315.88'
    'E16' = '  +0.57%  '
    'D17' = '0.813'
    'E17' = '  +0.59%  '
    'D18' = '42.920.21'
    'E18' = '  -0.84%  '
    'D19' = '6.22'
    'E19' = '  +0.60%  '
    'D20' = '0.0₃0912'
    'E20' = '  -2.55%  '
    'D21' = '11.75'
    'E21' = '  -4.60%  '
    'D22' = '67.86'
    'E22' = '  -0.31%  '
    'D23' = '236.99'
    'E23' = '  -1.90%  '
    'E24' = '  +0.63%  '
    'D25' = '2.55'
    'E25' = '  -2.36%  '
    'E26' = '  -0.04%  '
    'D27' = '25.43'
    'E27' = '  +3.05%  '
    'E28' = '  +1.07%  '
    'D29' = '35.02'
    'E29' = '  -4.73%  '
    'D30' = '9.38'
    'E30' = '  -2.55%  '
    'D31' = '163.32'
    'E31' = '  -3.03%  '
    'E32' = '  -0.02%  '
    'D33' = '5.12'
    'E33' = '  -2.94%  '
    'D34' = '17.57'
    'E34' = '  -0.46%  '
    'B35' = 'RenderToken'
    'C35' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D35' = '4.63'
    'E35' = '  +6.31%  '
    'B36' = 'WEMIXToken'
    'C36' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D36' = '2.47'
    'E36' = '  -1.59%  '
    'D37' = '0.0727'
    'E37' = '  -2.25%  '
    'B38' = 'LidoDAOToken'
    'C38' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D38' = '2.92'
    'E38' = '  -4.96%  '
    'B39' = 'ARBITRUM'
    'C39' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D39' = '1.85'
    'E39' = '  -2.40%  '
    'E40' = '  -3.05%  '
    'E41' = '  -2.06%  '
    'D42' = '2.44'
    'E42' = '  +5.36%  '
    'D43' = '2.002.65'
    'E43' = '  +1.69%  '
    'D44' = '0.0286'
    'E44' = '  -1.07%  '
    'D45' = '18.69'
    'E45' = '  -3.71%  '
    'D46' = '10.17'
    'E46' = '  +3.12%  '
    'D47' = '2.93'
    'E47' = '  -1.69%  '
    'D48' = '55.86'
    'E48' = '  +0.75%  '
    'E49' = '  -0.23%  '
    'D50' = '2.562.13'
    'E50' = '  +1.16%  '
    'D51' = '4.71'
    'E51' = '  +2.62%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to treat the assignment as literal text,
    # preventing numeric-looking strings (e.g. '306.50', '42.986.27') from being
    # reinterpreted/rounded as numbers.
    $range.Value = "'" + $updates[$cellRef]
    # Reset style to Normal so Excel's auto text-format style isn't left on the cell,
    # keeping cell styling identical to the original (unstyled) data cells.
    $range.Style = 'Normal'
}
